$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SkipBoss",    2,    960, 96, 100, 200, "win"),
    @("BonusPower",  0.02,   0, 17, 100,   0, "lose"),
    @("BonusPower",  2,    760, 96, 100, 200, "win"),
    @("BonusPower",  2,    980, 98, 100, 200, "win")
)

$startRow = 50
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
